$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1, matching the style of the other header cells (e.g. B1)
$ws.Range("H1").Value = "Save"
$ws.Range("B1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Fill the new Save column values for rows 2-7
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 0
